# Auto update Excel log
$wb = $excel.ActiveWorkbook

# --- Proximity sheet: append new ENTER/EXIT rows (11-15) ---
$proximity = $wb.Worksheets.Item("Proximity")

$proximityRows = @(
    @("2026-02-01", "13:36:44", "13:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "13:36:44", "13:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "13:36:56", "13:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door"),
    @("2026-02-01", "13:37:20", "13:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door"),
    @("2026-02-01", "13:37:32", "13:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)

$startRow = 11
$endRow = $startRow + $proximityRows.Count - 1
# Format the Date column as text first so the "YYYY-MM-DD" strings are kept
# as literal text instead of being auto-converted into date serial numbers.
$proximity.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $proximityRows.Count; $i++) {
    $row = $startRow + $i
    $data = $proximityRows[$i]
    $proximity.Cells.Item($row, 1).Value = $data[0]
    $proximity.Cells.Item($row, 2).Value = $data[1]
    $proximity.Cells.Item($row, 3).Value = $data[2]
    $proximity.Cells.Item($row, 4).Value = $data[3]
    $proximity.Cells.Item($row, 5).Value = $data[4]
    $proximity.Cells.Item($row, 6).Value = $data[5]
}

# --- Camera sheet: append new Image Captured rows (3-6) ---
$camera = $wb.Worksheets.Item("Camera")

$cameraRows = @(
    @("2026-02-01", "13:36:44", "13:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "13:36:56", "13:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "13:37:21", "13:00", "Living Room Main Door", "Image Captured", "Active"),
    @("2026-02-01", "13:37:32", "13:00", "Living Room Main Door", "Image Captured", "Active")
)

$startRow2 = 3
$endRow2 = $startRow2 + $cameraRows.Count - 1
$camera.Range("A$startRow2`:A$endRow2").NumberFormat = "@"

for ($i = 0; $i -lt $cameraRows.Count; $i++) {
    $row = $startRow2 + $i
    $data = $cameraRows[$i]
    $camera.Cells.Item($row, 1).Value = $data[0]
    $camera.Cells.Item($row, 2).Value = $data[1]
    $camera.Cells.Item($row, 3).Value = $data[2]
    $camera.Cells.Item($row, 4).Value = $data[3]
    $camera.Cells.Item($row, 5).Value = $data[4]
    $camera.Cells.Item($row, 6).Value = $data[5]
}
